$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 109 (old rows 109-112 shift down to 110-113,
# formulas/refs adjust automatically).
$ws.Rows.Item(109).Insert()

# Row 108: the end time was changed (E108).
$ws.Range("E108").Value = 0.60416666666666663

# Populate the newly inserted row 109 with a new working-hours entry.
$ws.Range("A109").Value = 2014
$ws.Range("B109").Value = 4
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 0.76041666666666663
$ws.Range("E109").Value = 0.79166666666666663
$ws.Range("F109").Formula = "=(E109-D109)*24*60"
$ws.Range("G109").Formula = "=F109/60"

# Fix the selection to match the post-edit cursor position.
$ws.Range("A110").Select()
